# lente contato inicio 1.1.9
# Update "Sttefani" (row 11) and "Suzana" (row 12) sales figures in the
# "VENDA EM BOLETOS - MATRIZ" table, then propagate the change through the
# dependent subtotal / summary rows and the pivot table at the bottom of
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Table 1 (rows 7-13): VENDA EM BOLETOS - MATRIZ -----------------------
# Row 11 - Sttefani
$ws.Range("C11").Value2 = 1950.0
$ws.Range("D11").Value2 = 600.0
$ws.Range("E11").Value2 = 1350.0
$ws.Range("F11").Value2 = 0.69230769230769
$ws.Range("G11").Value2 = 13.5

# Row 12 - Suzana
$ws.Range("C12").Value2 = 6000.0
$ws.Range("D12").Value2 = 4900.0
$ws.Range("E12").Value2 = 1100.0
$ws.Range("F12").Value2 = 0.18333333333333
$ws.Range("G12").Value2 = 11.0

# Row 13 - TOTAL of table 1
$ws.Range("C13").Value2 = 14828.0
$ws.Range("D13").Value2 = 9600.0
$ws.Range("E13").Value2 = 5228.0
$ws.Range("F13").Value2 = 0.35257620717561
$ws.Range("G13").Value2 = 52.28

# ---- Table 3 (rows 22-27): RESUMO GERAL VENDAS EM BOLETOS -----------------
# Row 24 - mirrors the MATRIZ total (row 13)
$ws.Range("C24").Value2 = 14828.0
$ws.Range("D24").Value2 = 9600.0
$ws.Range("E24").Value2 = 5228.0
$ws.Range("F24").Value2 = 0.35257620717561
$ws.Range("G24").Value2 = 52.28

# Row 27 - overall TOTAL (MATRIZ + PRESTIGIO + DAILY)
$ws.Range("C27").Value2 = 17968.0
$ws.Range("D27").Value2 = 11090.0
$ws.Range("E27").Value2 = 6878.0
$ws.Range("F27").Value2 = 0.38279162956367
$ws.Range("G27").Value2 = 68.78

# ---- Bottom pivot table (rows 30-41) --------------------------------------
# Row 34 - Sttefani
$ws.Range("D34").Value2 = 1950.0
$ws.Range("E34").Value2 = 600.0
$ws.Range("F34").Value2 = 1350.0

# Row 35 - Suzana
$ws.Range("D35").Value2 = 6000.0
$ws.Range("E35").Value2 = 4900.0
$ws.Range("F35").Value2 = 1100.0

# Recalculate the workbook so the SUM()/ratio formulas in rows 39-41 pick up
# the new inputs and their cached <v> results stay consistent.
$excel.Calculate()
